$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","F","G","H","I","L","N","O")

$newValues = @{
    2 = @(0.8829944718210641, 0.286624857758909, 0.01874274972951895, 0.459694673470004, 0.3054199823651373, 0.4648416140178, 0.443739043181381, 0.2998398206338777, 0.982931259631016, 1.468480130993214)
    3 = @(0.7867093547473019, 0.2794896162457832, 0.01662317875314301, 0.4549278335721496, 0.3019755812360785, 0.4667210908874182, 0.448295327536016, 0.2882069792364632, 0.9848650303441602, 1.464857129094469)
    4 = @(0.7275367405246413, 0.2751015698018762, 0.015314309747815, 0.4523365957836276, 0.3001249915484578, 0.4681185003265753, 0.4513918033336815, 0.2812338264995304, 0.9864337865403883, 1.463693244482101)
    5 = @(0.7034118145008961, 0.2733118480901169, 0.01477909462916926, 0.4513650190952987, 0.2994372104550607, 0.4687491895338098, 0.4527287319770963, 0.2784349095351786, 0.9871691956531521, 1.463485485407787)
    6 = @(0.6994052399025747, 0.273014577941737, 0.01469011233655237, 0.4512087848754476, 0.2993270090910656, 0.4688576141706378, 0.4529552611050747, 0.2779727330241428, 0.9872971221976172, 1.46346707994698)
    7 = @(0.7272114279873847, 0.2750774390315485, 0.01530709905913596, 0.452323151162858, 0.3001154473897714, 0.4681267580534225, 0.4514095296801521, 0.2811959063218836, 0.9864433150124512, 1.463689363623146)
    8 = @(0.8498074249825436, 0.2841661913589633, 0.01801348490715782, 0.4579813935825428, 0.3041774335312581, 0.4654391553667452, 0.4452479709922024, 0.2957936940852619, 0.9835189898908041, 1.467010669792714)
    9 = @(1.089730026220195, 0.3019259851008371, 0.02326048702741446, 0.4717424715242657, 0.3142457565532482, 0.4620991977128313, 0.4355400877213462, 0.32576244780401, 0.9808019967161385, 1.481950277817845)
    10 = @(1.265632098115134, 0.3149262661827947, 0.02707755183251948, 0.48348263877417, 0.3229342021733714, 0.4608215199038455, 0.4298601937169302, 0.3485988557300885, 0.9806349535199956, 1.498082287884984)
    11 = @(1.34556019050126, 0.3208281525361656, 0.02880556607559726, 0.4891786504845541, 0.3271692492411375, 0.4604955695827755, 0.4275926698953363, 0.3591655642312332, 0.9809540769983585, 1.506544962589004)
    12 = @(1.37581227144625, 0.323061131909185, 0.02945868581318223, 0.4913867333859159, 0.3288137286039046, 0.4604088316581567, 0.4267795816778985, 0.3631924912723719, 0.9811315561725849, 1.509911466073646)
    13 = @(1.369297640092555, 0.3225803093329205, 0.02931808058664842, 0.4909089089198488, 0.3284577459379108, 0.4604258805661345, 0.426952666407896, 0.3623240859733414, 0.9810908169048673, 1.50917922651567)
    14 = @(1.348049358067783, 0.3210119011126551, 0.02885932371967925, 0.4893592859818838, 0.3273037240131487, 0.4604876983204349, 0.4275248627463668, 0.3594963505663173, 0.9809675443171244, 1.506818681440365)
    15 = @(1.335032167580209, 0.3200509473201691, 0.02857815912122419, 0.4884167562126294, 0.3266021637275287, 0.4605303414513031, 0.4278812873023661, 0.3577676046411824, 0.9808994063749594, 1.505393866477618)
    16 = @(1.260406621264337, 0.314540304423673, 0.02696445021431515, 0.4831175440056086, 0.322663129173975, 0.4608479569348134, 0.4300147541390729, 0.3479118762905813, 0.9806220330088422, 1.49755187188174)
    17 = @(1.214601683914964, 0.3111564783728795, 0.0259723184647811, 0.4799576832514418, 0.3203191312111358, 0.4611081740975465, 0.4314046416370374, 0.3419113073623095, 0.9805529398458646, 1.493029141338752)
    18 = @(1.188247478258916, 0.3092090735230215, 0.02540088316975186, 0.4781736606676503, 0.3189975211220002, 0.4612818724427541, 0.432233833275351, 0.3384767282769303, 0.9805503997599345, 1.490533582962058)
    19 = @(1.179323009905318, 0.3085495301734227, 0.02520727085275354, 0.4775753647335108, 0.3185546108656183, 0.4613448110697078, 0.432519692691077, 0.3373167253602247, 0.9805559345844586, 1.489706793613465)
    20 = @(1.219478582744955, 0.3115168096265677, 0.02607801432916546, 0.48029059410581, 0.3205659004026273, 0.461077986806302, 0.431253604494426, 0.3425483412245285, 0.9805564461577205, 1.493499642551598)
    21 = @(1.354290914661021, 0.3214726348669785, 0.02899410571759375, 0.4898130600205306, 0.327641581115472, 0.4604685452503361, 0.4273555572211194, 0.3603262325277967, 0.9810022168668695, 1.50750763583676)
    22 = @(1.442310521034472, 0.3279679553783126, 0.03089268554126789, 0.4963345601698279, 0.3325035749329999, 0.460284103968462, 0.4250736394041823, 0.3720939554473688, 0.9816235785979188, 1.517606284361477)
    23 = @(1.395341490002579, 0.3245023917727963, 0.02988005333090626, 0.4928266354044055, 0.3298868561617923, 0.4603629803261811, 0.4262672018178684, 0.3657997134959885, 0.9812618080850797, 1.512130032557451)
    24 = @(1.217273799984525, 0.3113539100750984, 0.02603023248356351, 0.4801399834917746, 0.3204542551174399, 0.4610915594217602, 0.431321794550783, 0.3422602907764798, 0.9805547451189653, 1.49328660308791)
    25 = @(1.024884249432716, 0.2971292226984446, 0.0218476039027351, 0.4677338676302867, 0.3112959201759935, 0.462796146189703, 0.4379115911121083, 0.3175113828869058, 0.9812150612563499, 1.47700463823233)
}

foreach ($r in $newValues.Keys) {
    $rowVals = $newValues[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $rowVals[$i]
    }
}

Write-Host "Updated pl_mw values for case with 380 kV"
